$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 274.3
$ws.Cells.Item(18, 9).Value = 274.3
$ws.Cells.Item(18, 11).Value = 274.3
$ws.Cells.Item(18, 13).Value = 9.699999999999989
$ws.Cells.Item(43, 8).Value = 6873.75
$ws.Cells.Item(43, 9).Value = 6665
$ws.Cells.Item(43, 11).Value = 6665
$ws.Cells.Item(43, 13).Value = -6596
$ws.Cells.Item(51, 8).Value = 4999.5
$ws.Cells.Item(51, 9).Value = 5000
$ws.Cells.Item(51, 11).Value = 5000
$ws.Cells.Item(51, 13).Value = -4516
$ws.Cells.Item(127, 8).Value = 1666.3334
$ws.Cells.Item(127, 10).Value = 1666.3334
$ws.Cells.Item(127, 12).Value = 4999.0002
$ws.Cells.Item(127, 14).Value = -14919.0002
$ws.Cells.Item(133, 14).ClearContents()
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(137, 8).Value = 2118.7144
$ws.Cells.Item(137, 9).Value = 1603.909
$ws.Cells.Item(137, 11).Value = 4811.727000000001
$ws.Cells.Item(137, 13).Value = -2261.727000000001
$ws.Cells.Item(139, 14).ClearContents()
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2499.8462
$ws.Cells.Item(61, 9).Value = 2458.3333
$ws.Cells.Item(61, 10).Value = 2998
$ws.Cells.Item(61, 11).Value = 2458.3333
$ws.Cells.Item(61, 12).Value = 2998
$ws.Cells.Item(61, 13).Value = -2246.3333
$ws.Cells.Item(61, 14).Value = -3422
$ws.Cells.Item(74, 8).Value = 799.1667
$ws.Cells.Item(74, 9).Value = 759
$ws.Cells.Item(74, 10).Value = 1000
$ws.Cells.Item(74, 11).Value = 759
$ws.Cells.Item(74, 12).Value = 1000
$ws.Cells.Item(74, 13).Value = 115
$ws.Cells.Item(74, 14).Value = -2748
$ws.Cells.Item(77, 8).Value = 799.1667
$ws.Cells.Item(77, 9).Value = 759
$ws.Cells.Item(77, 10).Value = 1000
$ws.Cells.Item(77, 11).Value = 3795
$ws.Cells.Item(77, 12).Value = 5000
$ws.Cells.Item(77, 13).Value = 573
$ws.Cells.Item(77, 14).Value = -13736
$ws.Cells.Item(132, 8).Value = 2864.4583
$ws.Cells.Item(132, 9).Value = 1811.4286
$ws.Cells.Item(132, 11).Value = 5434.2858
$ws.Cells.Item(132, 13).Value = -2904.2858
$ws.Cells.Item(136, 8).Value = 2499.8462
$ws.Cells.Item(136, 9).Value = 2458.3333
$ws.Cells.Item(136, 10).Value = 2998
$ws.Cells.Item(136, 11).Value = 7374.999899999999
$ws.Cells.Item(136, 12).Value = 8994
$ws.Cells.Item(136, 13).Value = -4824.999899999999
$ws.Cells.Item(136, 14).Value = -14094
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(36, 14).ClearContents()
$ws.Cells.Item(36, 8).Value = 3000
$ws.Cells.Item(36, 9).Value = 3000
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 3000
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -2466
$ws.Cells.Item(44, 14).ClearContents()
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(86, 8).Value = 2599.625
$ws.Cells.Item(86, 9).Value = 2599.625
$ws.Cells.Item(86, 11).Value = 2599.625
$ws.Cells.Item(86, 13).Value = -1476.625
$ws.Cells.Item(89, 8).Value = 2599.625
$ws.Cells.Item(89, 9).Value = 2599.625
$ws.Cells.Item(89, 11).Value = 12998.125
$ws.Cells.Item(89, 13).Value = -7382.125
$ws.Cells.Item(134, 8).Value = 2522.5386
$ws.Cells.Item(134, 9).Value = 2199.8572
$ws.Cells.Item(134, 11).Value = 6599.571599999999
$ws.Cells.Item(134, 13).Value = -4064.571599999999
$ws.Cells.Item(138, 13).ClearContents()
$ws.Cells.Item(138, 8).Value = 124749.5
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 124749.5
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 124749.5
$ws.Cells.Item(138, 14).Value = -135029.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(47, 14).ClearContents()
$ws.Cells.Item(47, 8).Value = 1001
$ws.Cells.Item(47, 9).Value = 1001
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 11).Value = 1001
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 13).Value = -435
$ws.Cells.Item(55, 8).Value = 17495
$ws.Cells.Item(55, 9).Value = 10000
$ws.Cells.Item(55, 10).Value = 24990
$ws.Cells.Item(55, 11).Value = 10000
$ws.Cells.Item(55, 12).Value = 24990
$ws.Cells.Item(55, 13).Value = -9685
$ws.Cells.Item(55, 14).Value = -25620
$ws.Cells.Item(58, 8).Value = 2477.2354
$ws.Cells.Item(58, 9).Value = 2042.3334
$ws.Cells.Item(58, 11).Value = 2042.3334
$ws.Cells.Item(58, 13).Value = -1839.3334
$ws.Cells.Item(86, 8).Value = 19087.455
$ws.Cells.Item(86, 9).Value = 9609.5
$ws.Cells.Item(86, 10).Value = 26985.75
$ws.Cells.Item(86, 11).Value = 9609.5
$ws.Cells.Item(86, 12).Value = 26985.75
$ws.Cells.Item(86, 13).Value = -8486.5
$ws.Cells.Item(86, 14).Value = -29231.75
$ws.Cells.Item(89, 8).Value = 19087.455
$ws.Cells.Item(89, 9).Value = 9609.5
$ws.Cells.Item(89, 10).Value = 26985.75
$ws.Cells.Item(89, 11).Value = 48047.5
$ws.Cells.Item(89, 12).Value = 134928.75
$ws.Cells.Item(89, 13).Value = -42431.5
$ws.Cells.Item(89, 14).Value = -146160.75
$ws.Cells.Item(99, 8).Value = 5031.2856
$ws.Cells.Item(99, 9).Value = 5041.2
$ws.Cells.Item(99, 11).Value = 5041.2
$ws.Cells.Item(99, 13).Value = -3543.2
$ws.Cells.Item(126, 8).Value = 5031.2856
$ws.Cells.Item(126, 9).Value = 5041.2
$ws.Cells.Item(126, 11).Value = 15123.6
$ws.Cells.Item(126, 13).Value = -12653.6
$ws.Cells.Item(134, 8).Value = 3241.5
$ws.Cells.Item(134, 9).Value = 3241.5
$ws.Cells.Item(134, 11).Value = 9724.5
$ws.Cells.Item(134, 13).Value = -7189.5
$ws.Cells.Item(136, 8).Value = 2477.2354
$ws.Cells.Item(136, 9).Value = 2042.3334
$ws.Cells.Item(136, 11).Value = 6127.0002
$ws.Cells.Item(136, 13).Value = -3577.0002
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1666.6666
$ws.Cells.Item(68, 9).Value = 1500
$ws.Cells.Item(68, 10).Value = 2000
$ws.Cells.Item(68, 11).Value = 4500
$ws.Cells.Item(68, 12).Value = 6000
$ws.Cells.Item(68, 13).Value = -3689
$ws.Cells.Item(68, 14).Value = -7622
$ws.Cells.Item(71, 8).Value = 1666.6666
$ws.Cells.Item(71, 9).Value = 1500
$ws.Cells.Item(71, 10).Value = 2000
$ws.Cells.Item(71, 11).Value = 13500
$ws.Cells.Item(71, 12).Value = 18000
$ws.Cells.Item(71, 13).Value = -9444
$ws.Cells.Item(71, 14).Value = -26112
$ws.Cells.Item(137, 8).Value = 6031.636
$ws.Cells.Item(137, 9).Value = 6264.5713
$ws.Cells.Item(137, 10).Value = 5624
$ws.Cells.Item(137, 11).Value = 18793.7139
$ws.Cells.Item(137, 12).Value = 16872
$ws.Cells.Item(137, 13).Value = -13693.7139
$ws.Cells.Item(137, 14).Value = -27072
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(132, 8).Value = 2864.353
$ws.Cells.Item(132, 9).Value = 1966.7778
$ws.Cells.Item(132, 10).Value = 3874.125
$ws.Cells.Item(132, 11).Value = 5900.3334
$ws.Cells.Item(132, 12).Value = 11622.375
$ws.Cells.Item(132, 13).Value = -3370.3334
$ws.Cells.Item(132, 14).Value = -16682.375
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 13).ClearContents()
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(55, 8).Value = 500.33334
$ws.Cells.Item(55, 9).Value = 500.33334
$ws.Cells.Item(55, 11).Value = 500.33334
$ws.Cells.Item(55, 13).Value = -327.33334
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(132, 8).Value = 5939
$ws.Cells.Item(132, 9).Value = 5819
$ws.Cells.Item(132, 11).Value = 17457
$ws.Cells.Item(132, 13).Value = -14927
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 595.4545000000001
$ws.Cells.Item(107, 9).Value = 533.75
$ws.Cells.Item(107, 10).Value = 760
$ws.Cells.Item(107, 11).Value = 1601.25
$ws.Cells.Item(107, 12).Value = 2280
$ws.Cells.Item(107, 13).Value = 318.75
$ws.Cells.Item(107, 14).Value = -6120
$ws.Cells.Item(113, 8).Value = 1999
$ws.Cells.Item(113, 9).Value = 1999
$ws.Cells.Item(113, 11).Value = 5997
$ws.Cells.Item(113, 13).Value = -3827
$ws.Cells.Item(122, 8).Value = 3019.5833
$ws.Cells.Item(122, 9).Value = 3083.6
$ws.Cells.Item(122, 10).Value = 2699.5
$ws.Cells.Item(122, 11).Value = 9250.799999999999
$ws.Cells.Item(122, 12).Value = 8098.5
$ws.Cells.Item(122, 13).Value = -6800.799999999999
$ws.Cells.Item(122, 14).Value = -12998.5
$ws.Cells.Item(132, 8).Value = 4589.625
$ws.Cells.Item(132, 10).Value = 4838.9
$ws.Cells.Item(132, 12).Value = 14516.7
$ws.Cells.Item(132, 14).Value = -19576.7
$ws.Cells.Item(136, 8).Value = 2898.0435
$ws.Cells.Item(136, 9).Value = 2596.6155
$ws.Cells.Item(136, 11).Value = 7789.8465
$ws.Cells.Item(136, 13).Value = -5239.8465
